$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# The target cells store plain numeric-looking text ("49.23", "4.25", "53.47")
# as shared strings (not numbers). Assigning a bare numeric-looking string via
# .Value lets Excel auto-convert it to a real number, so we prefix with an
# apostrophe (exactly like typing it in the UI) to force text entry.
$ws.Range("B11").Value = "'49.23"
$ws.Range("C11").Value = "'4.25"
$ws.Range("D11").Value = "'53.47"
